# avgCsd14 and totalCsd14 over season
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing csd14 header to totalCsd14
$ws.Range("M1").Value = "totalCsd14"

# Insert a new column at N for avgCsd14 (shifts old N/O -> O/P)
$ws.Range("N1").EntireColumn.Insert()
$ws.Range("N1").Value = "avgCsd14"

# Populate avgCsd14 values (totalCsd14 / matches) for each player row
$ws.Range("N2").Value = 7
$ws.Range("N3").Value = 0.5
$ws.Range("N4").Value = 28.5
$ws.Range("N5").Value = -29.5
$ws.Range("N6").Value = -4
$ws.Range("N7").Value = -7
$ws.Range("N8").Value = -0.5
$ws.Range("N9").Value = -28.5
$ws.Range("N10").Value = 29.5
$ws.Range("N11").Value = 4
$ws.Range("N12").Value = -30.5
$ws.Range("N13").Value = 17.5
$ws.Range("N14").Value = -23.5
$ws.Range("N15").Value = 3
$ws.Range("N16").Value = -9.5
$ws.Range("N17").Value = 30.5
$ws.Range("N18").Value = -17.5
$ws.Range("N19").Value = 23.5
$ws.Range("N20").Value = -3
$ws.Range("N21").Value = 9.5
